$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = '_CROSSREF'

$ws.Range("E5").Value = '[Chaolin%Huang%NULL%0,                  Yeming%Wang%NULL%0,                  Xingwang%Li%NULL%0,                  Lili%Ren%NULL%0,                  Jianping%Zhao%NULL%0,                  Yi%Hu%NULL%0,                  Li%Zhang%NULL%0,                  Guohui%Fan%NULL%0,                  Jiuyang%Xu%NULL%0,                  Xiaoying%Gu%NULL%0,                  Zhenshun%Cheng%NULL%0,                  Ting%Yu%NULL%0,                  Jiaan%Xia%NULL%0,                  Yuan%Wei%NULL%0,                  Wenjuan%Wu%NULL%0,                  Xuelei%Xie%NULL%0,                  Wen%Yin%NULL%0,                  Hui%Li%NULL%0,                  Min%Liu%NULL%0,                  Yan%Xiao%NULL%0,                  Hong%Gao%NULL%0,                  Li%Guo%NULL%0,                  Jungang%Xie%NULL%0,                  Guangfa%Wang%NULL%0,                  Rongmeng%Jiang%NULL%0,                  Zhancheng%Gao%NULL%0,                  Qi%Jin%NULL%0,                  Jianwei%Wang%wangjw28@163.com%0,                  Bin%Cao%caobin_ben@163.com%0]'
$ws.Range("I5").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E6").Value = '[Yan%Deng%NULL%5,                  Wei%Liu%NULL%0,                  Kui%Liu%NULL%2,                  Yuan-Yuan%Fang%NULL%0,                  Jin%Shang%NULL%1,                  Ling%Zhou%NULL%1,                  Ke%Wang%NULL%0,                  Fan%Leng%NULL%1,                  Shuang%Wei%NULL%1,                  Lei%Chen%NULL%1,                  Hui-Guo%Liu%NULL%0,                  Pei-Fang%Wei%NULL%0,                  Pei-Fang%Wei%NULL%0]'
$ws.Range("I6").Value = '_PMC_CROSSREF'

$ws.Range("E7").Value = '[Yichun%Cheng%NULL%1,                  Ran%Luo%NULL%1,                  Kun%Wang%NULL%2,                  Meng%Zhang%NULL%1,                  Zhixiang%Wang%NULL%1,                  Lei%Dong%NULL%1,                  Junhua%Li%NULL%2,                  Ying%Yao%NULL%1,                  Shuwang%Ge%geshuwang@tjh.tjmu.edu.cn%1,                  Gang%Xu%xugang@tjh.tjmu.edu.cn%1]'
$ws.Range("I7").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E8").Value = '[Carly%Eastin%NULL%2,                  Travis%Eastin%NULL%1]'
$ws.Range("I8").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E9").Value = '[Aguilar%Medina Jos\u00e9 M.%coreGivesNoEmail%5,                Claro%Vald\u00e9s Ram\u00f3n%coreGivesNoEmail%2,                Obreg\u00f3n%Santos Angel G.%coreGivesNoEmail%2,                Proh\u00edas%Mart\u00ednez Juan%coreGivesNoEmail%2,                Rodr\u00edguez%Blanco Suilbert%coreGivesNoEmail%2]'
$ws.Range("I9").Value = '_CROSSREF'

$ws.Range("I10").Value = '_PMC_CROSSREF'

$ws.Range("E11").Value = '[Luwen%Wang%NULL%2,                  Xun%Li%NULL%2,                  Hui%Chen%NULL%1,                  Shaonan%Yan%NULL%1,                  Dong%Li%NULL%1,                  Yan%Li%NULL%1,                  Zuojiong%Gong%NULL%1]'
$ws.Range("I11").Value = '_PMC_CROSSREF'

$ws.Range("E12").Value = '[Xiao-Wei%Xu%NULL%0,                  Xiao-Xin%Wu%NULL%0,                  Xian-Gao%Jiang%NULL%0,                  Kai-Jin%Xu%NULL%0,                  Ling-Jun%Ying%NULL%0,                  Chun-Lian%Ma%NULL%0,                  Shi-Bo%Li%NULL%0,                  Hua-Ying%Wang%NULL%0,                  Sheng%Zhang%NULL%0,                  Hai-Nv%Gao%NULL%0,                  Ji-Fang%Sheng%NULL%0,                  Hong-Liu%Cai%NULL%0,                  Yun-Qing%Qiu%NULL%0,                  Lan-Juan%Li%NULL%0]'
$ws.Range("I12").Value = '_PMC_CROSSREF'

$ws.Range("E14").Value = '[Gemin%Zhang%NULL%3,                  Jie%Zhang%945128911@qq.com%1,                  Bowen%Wang%NULL%1,                  Xionglin%Zhu%NULL%1,                  Qiang%Wang%NULL%1,                  Shiming%Qiu%NULL%1]'
$ws.Range("I14").Value = '_PMC_Springer_CROSSREF'

$ws.Range("E15").Value = '[Xiaoli%Zhang%NULL%2,                  Huan%Cai%NULL%1,                  Jianhua%Hu%NULL%1,                  Jiangshan%Lian%NULL%1,                  Jueqing%Gu%NULL%1,                  Shanyan%Zhang%NULL%1,                  Chanyuan%Ye%NULL%1,                  Yingfeng%Lu%NULL%1,                  Ciliang%Jin%NULL%1,                  Guodong%Yu%NULL%1,                  Hongyu%Jia%NULL%1,                  Yimin%Zhang%NULL%1,                  Jifang%Sheng%jifang_sheng@zju.edu.cn%1,                  Lanjuan%Li%ljli@zju.edu.cn%0,                  Yida%Yang%yidayang65@zju.edu.cn%1]'
$ws.Range("I15").Value = '_PMC_elsevier_CROSSREF'

$ws.Range("E16").Value = '[Fei%Zhou%NULL%0,                  Ting%Yu%NULL%0,                  Ronghui%Du%NULL%2,                  Guohui%Fan%NULL%0,                  Ying%Liu%NULL%0,                  Zhibo%Liu%NULL%2,                  Jie%Xiang%NULL%0,                  Yeming%Wang%NULL%0,                  Bin%Song%NULL%2,                  Xiaoying%Gu%NULL%0,                  Lulu%Guan%NULL%2,                  Yuan%Wei%NULL%0,                  Hui%Li%NULL%0,                  Xudong%Wu%NULL%2,                  Jiuyang%Xu%NULL%0,                  Shengjin%Tu%NULL%2,                  Yi%Zhang%NULL%2,                  Hua%Chen%NULL%2,                  Bin%Cao%NULL%0]'
$ws.Range("I16").Value = '_PMC_elsevier_CROSSREF'
